# ans/result/1401CS42.xlsx — "Handles float input without breaking stuff"
#
# The marksheet's summary block (rows 10-12) gets real scored numbers
# instead of placeholder zeros, the "No./Marking/Total" labels pick up
# the title style used elsewhere in that block, and the per-question
# "Student Ans" column (A) gets filled in together with a pass/fail
# style (correctStyle / incorrectStyle) for every answered question.
# The 2nd and 3rd "Student Ans / Correct Ans" blocks (columns D-E from
# row 19 down, and the whole G:H block) are dropped entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary block (rows 10-12) ------------------------------------
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 14
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 10
$ws.Range("E10").Value = 28

$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 56
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "52/112"

# --- Drop the 2nd/3rd "Student Ans | Correct Ans" blocks -----------
# 3rd block (columns G:H) disappears completely.
$ws.Range("G1:H100").Clear()
# 2nd block (columns D:E) only survives for rows 16-18; everything
# below that reverts to blank.
$ws.Range("D19:E40").Clear()

# --- 2nd block answers that remain (rows 16-18) ---------------------
$ws.Range("D16").Value = "Option A"
$ws.Range("D16").Style = "correctStyle"

$ws.Range("D17").Value = "Option C"
$ws.Range("D17").Style = "correctStyle"

$ws.Range("D18").Value = "Option A"
$ws.Range("D18").Style = "incorrectStyle"

# --- 1st block: fill in "Student Ans" (column A) for every question
#     that was actually attempted, coloured by correctness -----------
$studentAnswers = @{
    16 = @("Option C", "incorrectStyle")
    18 = @("Option B", "correctStyle")
    22 = @("Option D", "correctStyle")
    23 = @("Option D", "correctStyle")
    24 = @("Option A", "correctStyle")
    26 = @("Option D", "incorrectStyle")
    27 = @("Option A", "correctStyle")
    29 = @("Option D", "correctStyle")
    30 = @("Option B", "correctStyle")
    31 = @("Option D", "correctStyle")
    32 = @("Option C", "correctStyle")
    35 = @("Option C", "incorrectStyle")
    38 = @("Option A", "correctStyle")
    39 = @("Option D", "correctStyle")
    40 = @("Option D", "correctStyle")
}

foreach ($row in $studentAnswers.Keys) {
    $answer = $studentAnswers[$row][0]
    $style = $studentAnswers[$row][1]
    $cell = $ws.Range("A$row")
    $cell.Value = $answer
    $cell.Style = $style
}
